# v1.3.0 fix grammar excel
# Remove the 6 bad/duplicate rows from the grammar sheet (column A),
# letting the remaining rows (already alphabetically sorted) shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @(
    "丁752丁丁福的",
    "人满开车打丁丁468丁",
    "{74}{74} 丁丁468丁 {55}{55}",
    "丁丁468丁{65}",
    "146 连起来",
    "丁丁468丁"
)

$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 1; $r--) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($targets -contains $val) {
        $ws.Rows($r).Delete()
    }
}

# Re-apply the AutoFilter over the shrunk range (it stays oversized on
# purpose, matching the original file's habit of extending well past the
# last data row: 91 -> 85, a flat -6 shift, same as the row count delta).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:A85").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name (used by the
# AutoFilter) in sync with the new range.
$sheetName = $ws.Name
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=" + $sheetName + "!`$A`$1:`$A`$85"
    }
}

# Restore the on-screen state after trimming the list: scrolled down with
# A31 selected.
[void]$ws.Activate()
[void]$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 22
